# Adds a new "comment" column header to every sheet in the workbook and
# switches the active/selected tab from "currency_conversions" to
# "sell_orders" (matching the upstream commit "forgot to commit pieces,
# commit local progress now").

$wb = $excel.ActiveWorkbook

# --- rsu: new header in G1 (row 1 currently spans A:F) ---
$ws = $wb.Worksheets.Item("rsu")
$ws.Cells.Item(1, 7).Value = "comment"
$ws.Cells.Item(1, 7).Font.Bold = $true

# --- espp: new header in G1 (row 1 currently spans A:F) ---
$ws = $wb.Worksheets.Item("espp")
$ws.Cells.Item(1, 7).Value = "comment"
$ws.Cells.Item(1, 7).Font.Bold = $true

# --- dividends: new header in F1 (row 1 currently spans A:E) ---
$ws = $wb.Worksheets.Item("dividends")
$ws.Cells.Item(1, 6).Value = "comment"
$ws.Cells.Item(1, 6).Font.Bold = $true

# --- buy_orders: G1 already exists (blank, pre-formatted) -> fill it in ---
$ws = $wb.Worksheets.Item("buy_orders")
$ws.Cells.Item(1, 7).Value = "comment"

# --- sell_orders: new header in G1 (row 1 currently spans A:F); this sheet
#     also becomes the active/selected tab ---
$ws = $wb.Worksheets.Item("sell_orders")
$ws.Cells.Item(1, 7).Value = "comment"
$ws.Cells.Item(1, 7).Font.Bold = $true
$ws.Activate()

# --- currency_conversions: new header in F1 (row 1 currently spans A:E);
#     loses the "selected" tab status to sell_orders above ---
$ws = $wb.Worksheets.Item("currency_conversions")
$ws.Cells.Item(1, 6).Value = "comment"
$ws.Cells.Item(1, 6).Font.Bold = $true
